$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header values P1 and Q1 (continuing sequence 0..13 -> 14,15), matching style of O1
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# For data rows 2..25: add new columns P (=2) and Q (=2),
# and swap values between I<->K and M<->O columns
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O

    $ws.Cells.Item($r, 9).Value = $kVal   # I becomes old K
    $ws.Cells.Item($r, 11).Value = $iVal  # K becomes old I
    $ws.Cells.Item($r, 13).Value = $oVal  # M becomes old O
    $ws.Cells.Item($r, 15).Value = $mVal  # O becomes old M

    $ws.Cells.Item($r, 16).Value = 2      # P
    $ws.Cells.Item($r, 17).Value = 2      # Q
}
